$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
# F1: "EQ" -> "EQ_num"
$ws.Range("F1").Value = "EQ_num"
# G1: new "Site_Name" column header
$ws.Range("G1").Value = "Site_Name"

# --- Data rows ---
# Add a 0-based index column (A), shorten the date strings from
# MM_DD_YYYY to MM_DD_YY, and replace the concatenated site-date-EQ
# string in column G with a friendly Site_Name lookup value.
$ws.Range("A2").Value = 0
$ws.Range("E2").Value = "08_22_23"
$ws.Range("G2").Value = "Northwell"
$ws.Range("A3").Value = 1
$ws.Range("E3").Value = "09_19_23"
$ws.Range("G3").Value = "Northwell"
$ws.Range("A4").Value = 2
$ws.Range("E4").Value = "09_07_23"
$ws.Range("G4").Value = "Circleville"
$ws.Range("A5").Value = 3
$ws.Range("E5").Value = "09_19_23"
$ws.Range("G5").Value = "Circleville"
$ws.Range("A6").Value = 4
$ws.Range("E6").Value = "08_23_23"
$ws.Range("G6").Value = "Houston"
$ws.Range("A7").Value = 5
$ws.Range("E7").Value = "09_19_23"
$ws.Range("G7").Value = "Houston"
$ws.Range("A8").Value = 6
$ws.Range("E8").Value = "08_23_23"
$ws.Range("G8").Value = "Hilliard"
$ws.Range("A9").Value = 7
$ws.Range("E9").Value = "09_07_23"
$ws.Range("G9").Value = "Grove city"
$ws.Range("A10").Value = 8
$ws.Range("E10").Value = "09_08_23"
$ws.Range("G10").Value = "Mentor OH"
$ws.Range("A11").Value = 9
$ws.Range("E11").Value = "08_18_23"
$ws.Range("G11").Value = "Encinogho"
$ws.Range("A12").Value = 10
$ws.Range("E12").Value = "08_22_23"
$ws.Range("G12").Value = "LA Site"

# --- Column width adjustments ---
# New column A gets a narrower custom width
$ws.Columns.Item(1).ColumnWidth = 6.5
# Column E widens slightly to fit the new date format
$ws.Columns.Item(5).ColumnWidth = 10.83

# --- Selection ---
$ws.Range("G17").Select() | Out-Null
